# The time entry export now only keeps the most recent two entries
# (idtimeentry 3 and 4), which previously lived in rows 6-7. Move that
# data up into rows 2-3 and drop the rest of the old rows (4-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the data that used to sit in rows 6-7 (idtimeentry 3 and 4) up to
# rows 2-3, preserving the original text formatting of the cells.
$ws.Range("A6:G7").Copy()
$ws.Range("A2").PasteSpecial(-4104)

# Remove the now-duplicated/obsolete rows 4-7; this shifts nothing further
# up (rows 2-3 already hold the final values) and shrinks the used range
# down to A1:G3.
$ws.Range("A4:G7").Delete()
